$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.18599999999999
$ws.Range("C3").Value = -11.4692
$ws.Range("C5").Value = -12.7052
$ws.Range("D7").Value = -7.093399999999992
$ws.Range("B9").Value = 8.393200000000007
$ws.Range("D9").Value = -8.749500000000006
$ws.Range("C11").Value = -13.2294
$ws.Range("C12").Value = -14.26690000000002
$ws.Range("B13").Value = 5.017300000000005
$ws.Range("B16").Value = 9.385000000000009
$ws.Range("B18").Value = 4.658400000000003
$ws.Range("B20").Value = 5.8859
$ws.Range("C21").Value = -12.72100000000001
$ws.Range("D21").Value = -8.212500000000004
